# Auto-generated edit script: applies Omega_Profits recompute deltas
# to the H:N "profit scenario" columns across multiple sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
# Row 69
$ws.Range("H69").Value = 19133.334
$ws.Range("J69").Value = 19142.857
$ws.Range("L69").Value = 57428.571
$ws.Range("N69").Value = -59176.571
# Row 72
$ws.Range("H72").Value = 19133.334
$ws.Range("J72").Value = 19142.857
$ws.Range("L72").Value = 172285.713
$ws.Range("N72").Value = -181021.713
# Row 80
$ws.Range("H80").Value = 538.6
$ws.Range("I80").Value = 398.66666
$ws.Range("K80").Value = 1195.99998
$ws.Range("M80").Value = -197.9999800000001
# Row 83
$ws.Range("H83").Value = 538.6
$ws.Range("I83").Value = 398.66666
$ws.Range("K83").Value = 3587.99994
$ws.Range("M83").Value = 1404.00006
# Row 112
$ws.Range("H112").Value = 2935.3845
$ws.Range("J112").Value = 3477.3
$ws.Range("L112").Value = 10431.9
$ws.Range("N112").Value = -12647.9
# Row 132
$ws.Range("H132").Value = 2301.8667
$ws.Range("I132").Value = 2346.1628
$ws.Range("K132").Value = 7038.4884
$ws.Range("M132").Value = -4508.4884
# Row 137
$ws.Range("H137").Value = 2506.2
$ws.Range("J137").Value = 2586.7273
$ws.Range("L137").Value = 7760.1819
$ws.Range("N137").Value = -12860.1819
# Row 138
$ws.Range("H138").Value = 3622.8545
$ws.Range("I138").Value = 1597.7646
$ws.Range("J138").Value = 4528.816
$ws.Range("K138").Value = 4793.293799999999
$ws.Range("L138").Value = 13586.448
$ws.Range("M138").Value = 346.7062000000005
$ws.Range("N138").Value = -23866.448
# Row 141
$ws.Range("H141").Value = 2947.9185
$ws.Range("I141").Value = 2829.5
$ws.Range("K141").Value = 8488.5
$ws.Range("M141").Value = -3308.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2969.0151
$ws.Range("I32").Value = 1315.1578
$ws.Range("K32").Value = 1315.1578
$ws.Range("M32").Value = -1028.1578
# Row 45
$ws.Range("H45").Value = 2676.7058
$ws.Range("I45").Value = 2636.9167
$ws.Range("K45").Value = 2636.9167
$ws.Range("M45").Value = -2259.9167
# Row 61
$ws.Range("H61").Value = 5868.0586
$ws.Range("I61").Value = 5197.143
$ws.Range("K61").Value = 5197.143
$ws.Range("M61").Value = -4985.143
# Row 74
$ws.Range("H74").Value = 1542.8334
$ws.Range("I74").Value = 1546.1428
$ws.Range("K74").Value = 1546.1428
$ws.Range("M74").Value = -672.1428000000001
# Row 77
$ws.Range("H77").Value = 1542.8334
$ws.Range("I77").Value = 1546.1428
$ws.Range("K77").Value = 7730.714
$ws.Range("M77").Value = -3362.714
# Row 88
$ws.Range("H88").Value = 1919.4445
$ws.Range("J88").Value = 2379.5
$ws.Range("L88").Value = 2379.5
$ws.Range("N88").Value = -3191.5
# Row 91
$ws.Range("H91").Value = 1919.4445
$ws.Range("J91").Value = 2379.5
$ws.Range("L91").Value = 2379.5
$ws.Range("N91").Value = -5187.5
# Row 132
$ws.Range("H132").Value = 3260.3333
$ws.Range("I132").Value = 3158.1177
$ws.Range("K132").Value = 9474.3531
$ws.Range("M132").Value = -6944.3531
# Row 136
$ws.Range("H136").Value = 5868.0586
$ws.Range("I136").Value = 5197.143
$ws.Range("K136").Value = 15591.429
$ws.Range("M136").Value = -13041.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2719.5386
$ws.Range("I20").Value = 2410.2354
$ws.Range("J20").Value = 3303.7778
$ws.Range("K20").Value = 2410.2354
$ws.Range("L20").Value = 3303.7778
$ws.Range("M20").Value = -2163.2354
$ws.Range("N20").Value = -3797.7778

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2870
$ws.Range("I31").Value = 3195.195
$ws.Range("J31").Value = 2489.0571
$ws.Range("K31").Value = 3195.195
$ws.Range("L31").Value = 2489.0571
$ws.Range("M31").Value = -2900.195
$ws.Range("N31").Value = -3079.0571
# Row 34
$ws.Range("H34").Value = 2870
$ws.Range("I34").Value = 3195.195
$ws.Range("J34").Value = 2489.0571
$ws.Range("K34").Value = 3195.195
$ws.Range("L34").Value = 2489.0571
$ws.Range("M34").Value = -2993.195
$ws.Range("N34").Value = -2893.0571
# Row 94
$ws.Range("H94").Value = 1700
$ws.Range("J94").Value = 1700
$ws.Range("L94").Value = 1700
$ws.Range("N94").Value = -2602
# Row 132
$ws.Range("H132").Value = 3615.3845
$ws.Range("I132").Value = 4100.0454
$ws.Range("J132").Value = 949.75
$ws.Range("K132").Value = 12300.1362
$ws.Range("L132").Value = 2849.25
$ws.Range("M132").Value = -9770.136200000001
$ws.Range("N132").Value = -7909.25
# Row 134
$ws.Range("H134").Value = 824.10254
$ws.Range("I134").Value = 814.5946
$ws.Range("K134").Value = 2443.7838
$ws.Range("M134").Value = 91.21619999999984

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Range("H50").Value = 999
$ws.Range("I50").Value = 999
$ws.Range("K50").Value = 2997
$ws.Range("M50").Value = -2516
# Row 53
$ws.Range("H53").Value = 999
$ws.Range("I53").Value = 999
$ws.Range("K53").Value = 2997
$ws.Range("M53").Value = -2516
# Row 74
$ws.Range("H74").Value = 9500
$ws.Range("J74").Value = 9500
$ws.Range("L74").Value = 28500
$ws.Range("N74").Value = -30622
# Row 77
$ws.Range("H77").Value = 9500
$ws.Range("J77").Value = 9500
$ws.Range("L77").Value = 85500
$ws.Range("N77").Value = -96108
# Row 94
$ws.Range("H94").Value = 14059.667
$ws.Range("I94").Value = 11468
$ws.Range("J94").Value = 15355.5
$ws.Range("K94").Value = 34404
$ws.Range("L94").Value = 46066.5
$ws.Range("M94").Value = -33728
$ws.Range("N94").Value = -47418.5
# Row 132
$ws.Range("H132").Value = 1792.6666
$ws.Range("I132").Value = 1633
$ws.Range("J132").Value = 1899.1111
$ws.Range("K132").Value = 14697
$ws.Range("L132").Value = 17091.9999
$ws.Range("M132").Value = -12167
$ws.Range("N132").Value = -22151.9999
# Row 139
$ws.Range("H139").Value = 5911.2666
$ws.Range("I139").Value = 3367.5
$ws.Range("J139").Value = 10998.8
$ws.Range("K139").Value = 10102.5
$ws.Range("L139").Value = 32996.39999999999
$ws.Range("M139").Value = -4962.5
$ws.Range("N139").Value = -43276.39999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 5193520.5
$ws.Range("J3").Value = 5214001
$ws.Range("L3").Value = 5214001
$ws.Range("N3").Value = -5214233
# Row 4
$ws.Range("H4").Value = 5004
$ws.Range("J4").Value = 5004
$ws.Range("L4").Value = 5004
$ws.Range("N4").Value = -5228
# Row 5
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 10000
$ws.Range("M5").Value = -9888
# Row 12
$ws.Range("H12").Value = 89141.71000000001
$ws.Range("J12").Value = 22502
$ws.Range("L12").Value = 22502
$ws.Range("N12").Value = -22782
# Row 14
$ws.Range("H14").Value = 202640.9
$ws.Range("J14").Value = 251727.5
$ws.Range("L14").Value = 251727.5
$ws.Range("N14").Value = -252063.5
# Row 21
$ws.Range("H21").Value = 25000
$ws.Range("J21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("N21").Value = -25346
# Row 22
$ws.Range("H22").Value = 2495
$ws.Range("I22").Value = 2495
$ws.Range("K22").Value = 2495
$ws.Range("M22").Value = -1966
# Row 30
$ws.Range("H30").Value = 25000
$ws.Range("J30").Value = 25000
$ws.Range("L30").Value = 25000
$ws.Range("N30").Value = -25210
# Row 54
$ws.Range("H54").Value = 49700
$ws.Range("J54").Value = 49700
$ws.Range("L54").Value = 49700
$ws.Range("N54").Value = -50480
# Row 126
$ws.Range("H126").Value = 6799.3335
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6799.3335
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 20398.0005
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -25338.0005
# Row 132
$ws.Range("H132").Value = 9083.862999999999
$ws.Range("I132").Value = 9087.857
$ws.Range("K132").Value = 27263.571
$ws.Range("M132").Value = -24733.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3368.4707
$ws.Range("I61").Value = 3360.0908
$ws.Range("J61").Value = 3383.8333
$ws.Range("K61").Value = 3360.0908
$ws.Range("L61").Value = 3383.8333
$ws.Range("M61").Value = -3158.0908
$ws.Range("N61").Value = -3787.8333
# Row 113
$ws.Range("H113").Value = 3368.4707
$ws.Range("I113").Value = 3360.0908
$ws.Range("J113").Value = 3383.8333
$ws.Range("K113").Value = 3360.0908
$ws.Range("L113").Value = 3383.8333
$ws.Range("M113").Value = -1190.0908
$ws.Range("N113").Value = -7723.8333
# Row 136
$ws.Range("H136").Value = 1475.9565
$ws.Range("I136").Value = 1186.75
$ws.Range("K136").Value = 3560.25
$ws.Range("M136").Value = -1010.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5805.4644
$ws.Range("I132").Value = 4738.364
$ws.Range("K132").Value = 14215.092
$ws.Range("M132").Value = -11685.092

Write-Output "Applied updates: sets + clears across all target sheets"